$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values are forced to Text to avoid Excel auto-converting
# numeric-looking strings (e.g. "509.26") into real numbers, which
# would lose trailing zeros / the textual representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.505.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.448.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.447.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("E13").Value = "  -7.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.877.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.509.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.436.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.75%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.530.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("E28").Value = "  -5.24%  "

$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("E35").Value = "  -0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.812"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.58%  "

$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "255.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.571"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0916"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("E51").Value = "  +1.12%  "
